# CharactersCreator and CharactersCreatorForm Ready
# Fill a new row (22) with placeholder "-" values for the skill columns,
# and move the selection/active cell to A2 (first sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns to fill on row 22 (CON + skill columns, skipping Deception/T)
$cols = @("I","P","Q","R","S","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG")
foreach ($col in $cols) {
    $ws.Range($col + "22").Value = "-"
}

# Update the selected/active cell shown when the sheet is reopened
$ws.Range("A2").Select()
